$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.632.17'
$ws.Range("E2").Value = '  +1.76%  '

$ws.Range("D3").Value = '3.460.14'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.18%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.461.12'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.566'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.38%  '

$ws.Range("E10").Value = '  +0.63%  '

$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.428'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.69%  '

$ws.Range("D13").Value = '4.058.82'
$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '

$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").Value = '65.520.64'
$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").Value = '3.483.92'
$ws.Range("E18").Value = '  +0.68%  '

$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '

$ws.Range("E25").Value = '  -1.42%  '

$ws.Range("E26").Value = '  +0.82%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("E28").Value = '  +1.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("E34").Value = '  +3.87%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  -3.81%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.17'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.892'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.25%  '

$ws.Range("E39").Value = '  -0.67%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.69%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.808.36'
$ws.Range("E42").Value = '  -0.88%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.38%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0307'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("E48").Value = '  +3.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '337.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.54%  '

$ws.Range("E50").Value = '  +1.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '32.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.50%  '
